$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement dataset for rows 2..24 (donor_id, category, survey_question, timestamp)
# Row 1 (headers) is left untouched.
$data = @(
    @(2, "0ce5dd49", "Outlining ideas or slides|Drafting full text|Summarising sources or meeting notes|Adjusting style for different audiences", "q07_write_subtasks", 45854.65247151545),
    @(3, "2c1001cb", "Outlining ideas or slides|Drafting full text|Proof-reading – tone adjustment|Summarising sources or meeting notes|Adjusting style for different audiences", "q07_write_subtasks", 45854.65247151545),
    @(4, "37cc37bf", "Outlining ideas or slides|Drafting full text|Proof-reading – tone adjustment|Summarising sources or meeting notes|Adjusting style for different audiences", "q07_write_subtasks", 45854.65247151545),
    @(5, "43faa0b9", "Outlining ideas or slides|Drafting full text|Proof-reading – tone adjustment|Summarising sources or meeting notes|Adjusting style for different audiences", "q07_write_subtasks", 45854.65247151545),
    @(6, "4abe3e88", "Drafting full text|Proof-reading – tone adjustment|Summarising sources or meeting notes|Adjusting style for different audiences", "q07_write_subtasks", 45854.65247151545),
    @(7, "50164f59", "I did not choose “Writing & professional communication”", "q07_write_subtasks", 45854.65247151545),
    @(8, "5cf70f79", "Outlining ideas or slides|Drafting full text|Proof-reading – tone adjustment|Summarising sources or meeting notes|Adjusting style for different audiences", "q07_write_subtasks", 45854.65247151545),
    @(9, "5da96769", "I did not choose “Writing & professional communication”", "q07_write_subtasks", 45854.65247151545),
    @(10, "6ca3e2f6", "Outlining ideas or slides|Drafting full text|Summarising sources or meeting notes|Adjusting style for different audiences", "q07_write_subtasks", 45854.65247151545),
    @(11, "790a4fcb", "Outlining ideas or slides|Drafting full text|Proof-reading – tone adjustment|Summarising sources or meeting notes|Adjusting style for different audiences", "q07_write_subtasks", 45854.65247151545),
    @(12, "802cc63a", "Outlining ideas or slides|Drafting full text|Proof-reading – tone adjustment|Summarising sources or meeting notes|Adjusting style for different audiences", "q07_write_subtasks", 45854.65247151545),
    @(13, "85c3ea4d", "I did not choose “Writing & professional communication”", "q07_write_subtasks", 45854.65247151545),
    @(14, "942dfafb", "I did not choose “Writing & professional communication”", "q07_write_subtasks", 45854.65247151545),
    @(15, "9bc6ba8c", "I did not choose “Writing & professional communication”", "q07_write_subtasks", 45854.65247151545),
    @(16, "a2d65af2", "Outlining ideas or slides|Drafting full text|Summarising sources or meeting notes|Adjusting style for different audiences", "q07_write_subtasks", 45854.65247151545),
    @(17, "a46f1771", "Outlining ideas or slides|Drafting full text|Summarising sources or meeting notes|Adjusting style for different audiences", "q07_write_subtasks", 45854.65247151545),
    @(18, "ad58f9da", "Outlining ideas or slides|Drafting full text|Summarising sources or meeting notes|Adjusting style for different audiences", "q07_write_subtasks", 45854.65247151545),
    @(19, "c7d9a301", "Outlining ideas or slides|Drafting full text|Proof-reading – tone adjustment|Summarising sources or meeting notes|Adjusting style for different audiences", "q07_write_subtasks", 45854.65247151545),
    @(20, "ce8732ff", "Outlining ideas or slides|Drafting full text|Proof-reading – tone adjustment|Summarising sources or meeting notes|Adjusting style for different audiences", "q07_write_subtasks", 45854.65247151545),
    @(21, "d6f1d567", "Outlining ideas or slides|Drafting full text|Proof-reading – tone adjustment|Summarising sources or meeting notes|Adjusting style for different audiences", "q07_write_subtasks", 45854.65247151545),
    @(22, "da9326c9", "Outlining ideas or slides|Drafting full text|Proof-reading – tone adjustment|Summarising sources or meeting notes|Adjusting style for different audiences", "q07_write_subtasks", 45854.65247151545),
    @(23, "e09ca7bf", "I did not choose “Writing & professional communication”", "q07_write_subtasks", 45854.65247151545),
    @(24, "ef53a641", "Outlining ideas or slides|Drafting full text|Summarising sources or meeting notes|Adjusting style for different audiences", "q07_write_subtasks", 45854.65247151545)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
